# Add a new "Player Info" worksheet as the first sheet in the workbook,
# populate it with player info for player 6921 (Luke Wood), and update the
# existing "ODI Batting" / "ODI Bowling" sheets so the match-card-link
# columns become a plain MATCH_CODE column (just the numeric match code).

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- Insert the new "Player Info" sheet before "ODI Batting" --------------
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# Inserting a sheet shifts sheet indices, so re-fetch the existing sheets by
# name now that "Player Info" has taken the first slot.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
$playerInfo.Range("A1:D1").Font.Bold = $true
$playerInfo.Range("A1:D1").HorizontalAlignment = -4108
$playerInfo.Range("A1:D1").VerticalAlignment = -4160
$playerInfo.Range("A1:D1").Borders.LineStyle = 1

$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "6921"
$playerInfo.Range("B2").Value = "Luke Wood"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Fast Medium"

# --- ODI Batting: MATCH_CARD_LINK -> MATCH_CODE (store just the code) ----
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4660"

# --- ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE (store just the code) ----
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4660"

$playerInfo.Range("A1").Select() | Out-Null
